# Auto-generated edit script: updates live crypto data values
# across sheets "Top 50 Cryptocurrencies", "Top 5 by Market Cap", and "Summary".
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Top 50 Cryptocurrencies" ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2, 3).Value = 98821 ; $ws1.Cells.Item(2, 4).Value = 1955065374675 ; $ws1.Cells.Item(2, 5).Value = 107185165245 ; $ws1.Cells.Item(2, 6).Value = 1.52129
$ws1.Cells.Item(3, 3).Value = 3367.93 ; $ws1.Cells.Item(3, 4).Value = 405569114425 ; $ws1.Cells.Item(3, 5).Value = 56041266202 ; $ws1.Cells.Item(3, 6).Value = 7.93654
$ws1.Cells.Item(4, 4).Value = 130803059570 ; $ws1.Cells.Item(4, 5).Value = 112196096160 ; $ws1.Cells.Item(4, 6).Value = 0.05625
$ws1.Cells.Item(5, 3).Value = 259.04 ; $ws1.Cells.Item(5, 4).Value = 122966673076 ; $ws1.Cells.Item(5, 5).Value = 14806263273 ; $ws1.Cells.Item(5, 6).Value = 8.57286
$ws1.Cells.Item(6, 3).Value = 629.9299999999999 ; $ws1.Cells.Item(6, 4).Value = 91896601685 ; $ws1.Cells.Item(6, 5).Value = 2508195564 ; $ws1.Cells.Item(6, 6).Value = 3.24204
$ws1.Cells.Item(7, 3).Value = 1.4 ; $ws1.Cells.Item(7, 4).Value = 79973167811 ; $ws1.Cells.Item(7, 5).Value = 18537037700 ; $ws1.Cells.Item(7, 6).Value = 26.17071
$ws1.Cells.Item(8, 3).Value = 0.394029 ; $ws1.Cells.Item(8, 4).Value = 57858872160 ; $ws1.Cells.Item(8, 5).Value = 9755401127 ; $ws1.Cells.Item(8, 6).Value = 2.45483
$ws1.Cells.Item(9, 3).Value = 0.9998 ; $ws1.Cells.Item(9, 4).Value = 38307987931 ; $ws1.Cells.Item(9, 5).Value = 15326576712 ; $ws1.Cells.Item(9, 6).Value = -0.04826
$ws1.Cells.Item(10, 3).Value = 3367.81 ; $ws1.Cells.Item(10, 4).Value = 33004494921 ; $ws1.Cells.Item(10, 5).Value = 138827377 ; $ws1.Cells.Item(10, 6).Value = 7.92654
$ws1.Cells.Item(11, 3).Value = 0.880489 ; $ws1.Cells.Item(11, 4).Value = 31509087095 ; $ws1.Cells.Item(11, 5).Value = 3762401528 ; $ws1.Cells.Item(11, 6).Value = 11.75574
$ws1.Cells.Item(12, 3).Value = 0.199221 ; $ws1.Cells.Item(12, 4).Value = 17202885060 ; $ws1.Cells.Item(12, 5).Value = 1073391665 ; $ws1.Cells.Item(12, 6).Value = 0.92219
$ws1.Cells.Item(13, 3).Value = 36.16 ; $ws1.Cells.Item(13, 4).Value = 14793159228 ; $ws1.Cells.Item(13, 5).Value = 1040535609 ; $ws1.Cells.Item(13, 6).Value = 6.46775
$ws1.Cells.Item(14, 4).Value = 14652797101 ; $ws1.Cells.Item(14, 5).Value = 1598117103 ; $ws1.Cells.Item(14, 6).Value = 3.12022
$ws1.Cells.Item(15, 3).Value = 4000.6 ; $ws1.Cells.Item(15, 4).Value = 14442456649 ; $ws1.Cells.Item(15, 5).Value = 172086574 ; $ws1.Cells.Item(15, 6).Value = 8.18927
$ws1.Cells.Item(16, 3).Value = 98638 ; $ws1.Cells.Item(16, 4).Value = 14408704518 ; $ws1.Cells.Item(16, 5).Value = 826177901 ; $ws1.Cells.Item(16, 6).Value = 1.6381
$ws1.Cells.Item(17, 3).Value = 5.52 ; $ws1.Cells.Item(17, 4).Value = 14069277184 ; $ws1.Cells.Item(17, 5).Value = 625981057 ; $ws1.Cells.Item(17, 6).Value = 1.96154
$ws1.Cells.Item(18, 3).Value = 3.59 ; $ws1.Cells.Item(18, 4).Value = 10230076114 ; $ws1.Cells.Item(18, 5).Value = 2195757402 ; $ws1.Cells.Item(18, 6).Value = 0.85599
$ws1.Cells.Item(19, 3).Value = 493.97 ; $ws1.Cells.Item(19, 4).Value = 9775691671 ; $ws1.Cells.Item(19, 5).Value = 1709315082 ; $ws1.Cells.Item(19, 6).Value = -6.39473
$ws1.Cells.Item(20, 3).Value = 3365.9 ; $ws1.Cells.Item(20, 4).Value = 9644552154 ; $ws1.Cells.Item(20, 5).Value = 1442232704 ; $ws1.Cells.Item(20, 6).Value = 8.026949999999999
$ws1.Cells.Item(21, 3).Value = 15.27 ; $ws1.Cells.Item(21, 4).Value = 9572585986 ; $ws1.Cells.Item(21, 5).Value = 1260125079 ; $ws1.Cells.Item(21, 6).Value = 4.8129
$ws1.Cells.Item(22, 3).Value = 0.00002131 ; $ws1.Cells.Item(22, 4).Value = 8965828127 ; $ws1.Cells.Item(22, 5).Value = 6726966462 ; $ws1.Cells.Item(22, 6).Value = 9.410690000000001
$ws1.Cells.Item(23, 3).Value = 6.21 ; $ws1.Cells.Item(23, 4).Value = 8949858104 ; $ws1.Cells.Item(23, 5).Value = 844734603 ; $ws1.Cells.Item(23, 6).Value = 9.40584
$ws1.Cells.Item(24, 3).Value = 0.286915 ; $ws1.Cells.Item(24, 4).Value = 8617828912 ; $ws1.Cells.Item(24, 5).Value = 2346286480 ; $ws1.Cells.Item(24, 6).Value = 20.77031
$ws1.Cells.Item(25, 3).Value = 8.800000000000001 ; $ws1.Cells.Item(25, 4).Value = 8133488447 ; $ws1.Cells.Item(25, 5).Value = 3391395 ; $ws1.Cells.Item(25, 6).Value = 3.53993
$ws1.Cells.Item(26, 3).Value = 5.76 ; $ws1.Cells.Item(26, 4).Value = 7015382282 ; $ws1.Cells.Item(26, 5).Value = 1003512868 ; $ws1.Cells.Item(26, 6).Value = 4.87787
$ws1.Cells.Item(27, 3).Value = 90.18000000000001 ; $ws1.Cells.Item(27, 4).Value = 6783548458 ; $ws1.Cells.Item(27, 5).Value = 1315864262 ; $ws1.Cells.Item(27, 6).Value = 1.30091
$ws1.Cells.Item(28, 3).Value = 12.02 ; $ws1.Cells.Item(28, 4).Value = 6422121172 ; $ws1.Cells.Item(28, 5).Value = 848170177 ; $ws1.Cells.Item(28, 6).Value = 3.68462
$ws1.Cells.Item(29, 3).Value = 3546.31 ; $ws1.Cells.Item(29, 4).Value = 6169611515 ; $ws1.Cells.Item(29, 5).Value = 103651261 ; $ws1.Cells.Item(29, 6).Value = 8.18763
$ws1.Cells.Item(30, 3).Value = 9.32 ; $ws1.Cells.Item(30, 4).Value = 5595254907 ; $ws1.Cells.Item(30, 5).Value = 872429790 ; $ws1.Cells.Item(30, 6).Value = 5.52257
$ws1.Cells.Item(31, 3).Value = 0.138422 ; $ws1.Cells.Item(31, 4).Value = 5288274559 ; $ws1.Cells.Item(31, 5).Value = 964046702 ; $ws1.Cells.Item(31, 6).Value = 10.79753
$ws1.Cells.Item(32, 3).Value = 0.9990520000000001 ; $ws1.Cells.Item(32, 4).Value = 5226406988 ; $ws1.Cells.Item(32, 5).Value = 16006065 ; $ws1.Cells.Item(32, 6).Value = -0.36705
$ws1.Cells.Item(33, 3).Value = 0.187995 ; $ws1.Cells.Item(33, 4).Value = 5141784943 ; $ws1.Cells.Item(33, 5).Value = 150634106 ; $ws1.Cells.Item(33, 6).Value = 7.82868
$ws1.Cells.Item(34, 3).Value = 9.6 ; $ws1.Cells.Item(34, 4).Value = 4555261310 ; $ws1.Cells.Item(34, 5).Value = 273488349 ; $ws1.Cells.Item(34, 6).Value = 6.60921
$ws1.Cells.Item(35, 3).Value = 27.94 ; $ws1.Cells.Item(35, 4).Value = 4180801808 ; $ws1.Cells.Item(35, 5).Value = 861244532 ; $ws1.Cells.Item(35, 6).Value = 4.70517
$ws1.Cells.Item(36, 3).Value = 0.00005242 ; $ws1.Cells.Item(36, 4).Value = 3950095141 ; $ws1.Cells.Item(36, 5).Value = 1601466165 ; $ws1.Cells.Item(36, 6).Value = -0.10235
$ws1.Cells.Item(37, 3).Value = 0.151879 ; $ws1.Cells.Item(37, 4).Value = 3832713780 ; $ws1.Cells.Item(37, 5).Value = 150125959 ; $ws1.Cells.Item(37, 6).Value = 0.45116
$ws1.Cells.Item(38, 3).Value = 7.34 ; $ws1.Cells.Item(38, 4).Value = 3801653465 ; $ws1.Cells.Item(38, 5).Value = 429954829 ; $ws1.Cells.Item(38, 6).Value = 0.74086
$ws1.Cells.Item(39, 3).Value = 0.465237 ; $ws1.Cells.Item(39, 4).Value = 3705435184 ; $ws1.Cells.Item(39, 5).Value = 501118230 ; $ws1.Cells.Item(39, 6).Value = 6.34705
$ws1.Cells.Item(40, 3).Value = 500.21 ; $ws1.Cells.Item(40, 4).Value = 3696581786 ; $ws1.Cells.Item(40, 5).Value = 280047488 ; $ws1.Cells.Item(40, 6).Value = 3.58771
$ws1.Cells.Item(41, 3).Value = 1.002 ; $ws1.Cells.Item(41, 4).Value = 3688254516 ; $ws1.Cells.Item(41, 5).Value = 222230206 ; $ws1.Cells.Item(41, 6).Value = -0.18498
$ws1.Cells.Item(42, 3).Value = 24.78 ; $ws1.Cells.Item(42, 4).Value = 3570818581 ; $ws1.Cells.Item(42, 5).Value = 32982943 ; $ws1.Cells.Item(42, 6).Value = 2.58192
$ws1.Cells.Item(43, 3).Value = 0.999769 ; $ws1.Cells.Item(43, 4).Value = 3441159133 ; $ws1.Cells.Item(43, 5).Value = 151317872 ; $ws1.Cells.Item(43, 6).Value = -0.09085
$ws1.Cells.Item(44, 4).Value = 3439796552 ; $ws1.Cells.Item(44, 5).Value = 304975628 ; $ws1.Cells.Item(44, 6).Value = 5.48871
$ws1.Cells.Item(45, 3).Value = 3.34 ; $ws1.Cells.Item(45, 4).Value = 3349732168 ; $ws1.Cells.Item(45, 5).Value = 1291008898 ; $ws1.Cells.Item(45, 6).Value = 5.12188
$ws1.Cells.Item(46, 4).Value = 3329445059 ; $ws1.Cells.Item(46, 5).Value = 482295249 ; $ws1.Cells.Item(46, 6).Value = 3.22596
$ws1.Cells.Item(47, 3).Value = 0.783631 ; $ws1.Cells.Item(47, 4).Value = 3213465867 ; $ws1.Cells.Item(47, 5).Value = 1670862527 ; $ws1.Cells.Item(47, 6).Value = 13.38334
$ws1.Cells.Item(48, 3).Value = 160.73 ; $ws1.Cells.Item(48, 4).Value = 2965520995 ; $ws1.Cells.Item(48, 5).Value = 85739377 ; $ws1.Cells.Item(48, 6).Value = -0.53269
$ws1.Cells.Item(49, 4).Value = 2935957058 ; $ws1.Cells.Item(49, 5).Value = 349512301 ; $ws1.Cells.Item(49, 6).Value = 1.94279
$ws1.Cells.Item(50, 3).Value = 4.67 ; $ws1.Cells.Item(50, 4).Value = 2807537544 ; $ws1.Cells.Item(50, 5).Value = 567842261 ; $ws1.Cells.Item(50, 6).Value = 5.23522
$ws1.Cells.Item(51, 3).Value = 0.832329 ; $ws1.Cells.Item(51, 4).Value = 2803686073 ; $ws1.Cells.Item(51, 5).Value = 186631406 ; $ws1.Cells.Item(51, 6).Value = 13.99973

# --- Sheet 2: "Top 5 by Market Cap" ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(2, 2).Value = 1955065374675
$ws2.Cells.Item(3, 2).Value = 405569114425
$ws2.Cells.Item(4, 2).Value = 130803059570
$ws2.Cells.Item(5, 2).Value = 122966673076
$ws2.Cells.Item(6, 2).Value = 91896601685

# --- Sheet 3: "Summary" ---
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(2, 2).Value = "'$4348.79"
$ws3.Cells.Item(3, 2).Value = "XRP (26.17%)"
$ws3.Cells.Item(4, 2).Value = "Bitcoin Cash (-6.39%)"
